# Updated cryptos list on Fri Jan  5 05:48:52 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.779.48"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").Value = "'2.247.46"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'322.99"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").Value = "'101.86"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("D7").Value = "'0.580"
$ws.Range("E7").Value = "  -1.15%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("D10").Value = "'37.06"
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("D11").Value = "'0.0832"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("E13").Value = "  -2.62%  "
$ws.Range("D14").Value = "'2.589.77"
$ws.Range("E14").Value = "  +0.25%  "
$ws.Range("D15").Value = "'0.856"
$ws.Range("E15").Value = "  -1.55%  "
$ws.Range("D16").Value = "'14.20"
$ws.Range("E16").Value = "  -1.69%  "
$ws.Range("D17").Value = "'2.254.74"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "'43.703.28"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("D19").Value = "'13.59"
$ws.Range("E19").Value = "  -6.33%  "
$ws.Range("E20").Value = "  +1.81%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "'65.54"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("E23").Value = "  -1.05%  "
$ws.Range("D24").Value = "'236.09"
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").Value = "'10.17"
$ws.Range("E28").Value = "  -2.61%  "
$ws.Range("D29").Value = "'37.04"
$ws.Range("E29").Value = "  +3.58%  "
$ws.Range("D30").Value = "'6.31"
$ws.Range("E30").Value = "  -1.73%  "
$ws.Range("D31").Value = "'160.10"
$ws.Range("E31").Value = "  +3.86%  "
$ws.Range("D32").Value = "'20.17"
$ws.Range("E32").Value = "  -1.62%  "
$ws.Range("E33").Value = "  -3.09%  "
$ws.Range("E34").Value = "  -1.51%  "
$ws.Range("E35").Value = "  +9.62%  "
$ws.Range("D36").Value = "'3.06"
$ws.Range("E36").Value = "  -4.01%  "
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("E38").Value = "  -2.98%  "
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("E40").Value = "  -4.79%  "
$ws.Range("D41").Value = "'15.87"
$ws.Range("E41").Value = "  +20.05%  "
$ws.Range("E42").Value = "  -2.43%  "
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("D44").Value = "'1.811.06"
$ws.Range("E44").Value = "  +1.56%  "
$ws.Range("D45").Value = "'76.25"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("E46").Value = "  -3.23%  "
$ws.Range("D47").Value = "'82.44"
$ws.Range("E47").Value = "  -5.07%  "
$ws.Range("D48").Value = "'5.21"
$ws.Range("E48").Value = "  -2.55%  "
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "'58.61"
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'1.68"
$ws.Range("E50").Value = "  +4.63%  "
$ws.Range("D51").Value = "'103.45"
$ws.Range("E51").Value = "  -0.60%  "
